# Auto-generated: apply Leve profit recalculation updates across multiple worksheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 13421.2
$ws.Range("I111").Value = 27139.5
$ws.Range("K111").Value = 81418.5
$ws.Range("M111").Value = -78351.5
$ws.Range("H113").Value = 252326.25
$ws.Range("I113").Value = 252326.25
$ws.Range("K113").Value = 252326.25
$ws.Range("M113").Value = -249072.25
$ws.Range("H132").Value = 5439699.5
$ws.Range("I132").Value = 5957682.5
$ws.Range("K132").Value = 17873047.5
$ws.Range("M132").Value = -17870517.5
$ws.Range("H137").Value = 1807.3334
$ws.Range("I137").Value = 1436.1428
$ws.Range("J137").Value = 3106.5
$ws.Range("K137").Value = 4308.428400000001
$ws.Range("L137").Value = 9319.5
$ws.Range("M137").Value = -1758.428400000001
$ws.Range("N137").Value = -14419.5
$ws.Range("H138").Value = 4431.604
$ws.Range("I138").Value = 2398.3076
$ws.Range("J138").Value = 5092.425
$ws.Range("K138").Value = 7194.9228
$ws.Range("L138").Value = 15277.275
$ws.Range("M138").Value = -2054.9228
$ws.Range("N138").Value = -25557.275

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 54128.05
$ws.Range("I2").Value = 1353.3334
$ws.Range("J2").Value = 144599
$ws.Range("K2").Value = 1353.3334
$ws.Range("L2").Value = 144599
$ws.Range("M2").Value = -1240.3334
$ws.Range("N2").Value = -144825
$ws.Range("H45").Value = 1937.625
$ws.Range("I45").Value = 2123
$ws.Range("K45").Value = 2123
$ws.Range("M45").Value = -1746
$ws.Range("H61").Value = 2443.5557
$ws.Range("I61").Value = 2124
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 2124
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -1912
$ws.Range("N61").Value = -5424
$ws.Range("H110").Value = 50055540
$ws.Range("I110").Value = 50055540
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 50055540
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -50053495
$ws.Range("N110").ClearContents()
$ws.Range("H116").Value = 54128.05
$ws.Range("I116").Value = 1353.3334
$ws.Range("J116").Value = 144599
$ws.Range("K116").Value = 1353.3334
$ws.Range("L116").Value = 144599
$ws.Range("M116").Value = 940.6666
$ws.Range("N116").Value = -149187
$ws.Range("H136").Value = 2443.5557
$ws.Range("I136").Value = 2124
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 6372
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -3822
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 54128.05
$ws.Range("I3").Value = 1353.3334
$ws.Range("J3").Value = 144599
$ws.Range("K3").Value = 1353.3334
$ws.Range("L3").Value = 144599
$ws.Range("M3").Value = -1239.3334
$ws.Range("N3").Value = -144827
$ws.Range("H99").Value = 2416
$ws.Range("I99").Value = 2293.3333
$ws.Range("J99").Value = 2600
$ws.Range("K99").Value = 2293.3333
$ws.Range("L99").Value = 2600
$ws.Range("M99").Value = -795.3332999999998
$ws.Range("N99").Value = -5596
$ws.Range("H107").Value = 55626076
$ws.Range("I107").Value = 83436090
$ws.Range("K107").Value = 83436090
$ws.Range("M107").Value = -83434170

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 35991.332
$ws.Range("I99").Value = 4980
$ws.Range("J99").Value = 51497
$ws.Range("K99").Value = 4980
$ws.Range("L99").Value = 51497
$ws.Range("M99").Value = -3482
$ws.Range("N99").Value = -54493
$ws.Range("H126").Value = 35991.332
$ws.Range("I126").Value = 4980
$ws.Range("J126").Value = 51497
$ws.Range("K126").Value = 14940
$ws.Range("L126").Value = 154491
$ws.Range("M126").Value = -12470
$ws.Range("N126").Value = -159431

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 1492.3334
$ws.Range("I116").Value = 399.33334
$ws.Range("K116").Value = 1198.00002
$ws.Range("M116").Value = 2243.99998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3508.75
$ws.Range("I80").Value = 3492.5
$ws.Range("K80").Value = 3492.5
$ws.Range("M80").Value = -2494.5
$ws.Range("H83").Value = 3508.75
$ws.Range("I83").Value = 3492.5
$ws.Range("K83").Value = 17462.5
$ws.Range("M83").Value = -12470.5
$ws.Range("H113").Value = 1508.7273
$ws.Range("I113").Value = 1266.5
$ws.Range("J113").Value = 1799.4
$ws.Range("K113").Value = 1266.5
$ws.Range("L113").Value = 1799.4
$ws.Range("M113").Value = 903.5
$ws.Range("N113").Value = -6139.4
$ws.Range("H122").Value = 6663.375
$ws.Range("I122").Value = 5801.4
$ws.Range("J122").Value = 8100
$ws.Range("K122").Value = 17404.2
$ws.Range("L122").Value = 24300
$ws.Range("M122").Value = -14954.2
$ws.Range("N122").Value = -29200

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2233.35
$ws.Range("I7").Value = 1813.3572
$ws.Range("J7").Value = 3213.3333
$ws.Range("K7").Value = 1813.3572
$ws.Range("L7").Value = 3213.3333
$ws.Range("M7").Value = -1701.3572
$ws.Range("N7").Value = -3437.3333
$ws.Range("H126").Value = 2233.35
$ws.Range("I126").Value = 1813.3572
$ws.Range("J126").Value = 3213.3333
$ws.Range("K126").Value = 5440.071599999999
$ws.Range("L126").Value = 9639.999899999999
$ws.Range("M126").Value = -2970.071599999999
$ws.Range("N126").Value = -14579.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H87").Value = 19000
$ws.Range("J87").Value = 19000
$ws.Range("L87").Value = 19000
$ws.Range("N87").Value = -21496
$ws.Range("H90").Value = 19000
$ws.Range("J90").Value = 19000
$ws.Range("L90").Value = 57000
$ws.Range("N90").Value = -69480
$ws.Range("H100").Value = 125749.875
$ws.Range("I100").Value = 250549.75
$ws.Range("J100").Value = 950
$ws.Range("K100").Value = 501099.5
$ws.Range("L100").Value = 1900
$ws.Range("M100").Value = -500558.5
$ws.Range("N100").Value = -2982
$ws.Range("H113").Value = 484.25
$ws.Range("I113").Value = 366.33334
$ws.Range("K113").Value = 1099.00002
$ws.Range("M113").Value = 1070.99998
$ws.Range("H132").Value = 3107.8845
$ws.Range("I132").Value = 3055.35
$ws.Range("K132").Value = 9166.049999999999
$ws.Range("M132").Value = -6636.049999999999
